# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" columns for the
# 3bd7285c-da35-40db-bb99-7ef46d13247b handback row (row 6) on both the
# zh-cn and de-de sheets, because the handback file that came in was not
# based on the latest handoff version.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a6c36a74080edce79dc3f50d05c4906488bf2c0a/e2e/3bd7285c-da35-40db-bb99-7ef46d13247b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50fd26afdce35f10eb6f2fb02beb989cd0b1dabf/e2e/3bd7285c-da35-40db-bb99-7ef46d13247b.md."
$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50fd26afdce35f10eb6f2fb02beb989cd0b1dabf/e2e/3bd7285c-da35-40db-bb99-7ef46d13247b.md"
$mdDisplay = "3bd7285c-da35-40db-bb99-7ef46d13247b.md"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets("zh-cn")

# I6 -- "Latest Target File": link to the latest version of the handback markdown
$ws.Hyperlinks.Add($ws.Range("I6"), $latestMdUrl, "", "", $mdDisplay)

# J6 -- "Latest Handback File"
$ws.Range("J6").Value = "3bd7285c-da35-40db-bb99-7ef46d13247b.2928f6dea572a6e19d0edfc03e4fde9c41c5d031.zh-cn.xlf"

# K6 -- "Latest Handback DateTime"
$ws.Range("K6").Value = "2016-09-01 22:47:24"

# P6 -- "Error Detail"
$ws.Range("P6").Value = $errorDetail

# Widen the Error Detail column so the message is readable
$ws.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$ws2 = $wb.Worksheets("de-de")

# I6 -- "Latest Target File"
$ws2.Hyperlinks.Add($ws2.Range("I6"), $latestMdUrl, "", "", $mdDisplay)

# J6 -- "Latest Handback File"
$ws2.Range("J6").Value = "3bd7285c-da35-40db-bb99-7ef46d13247b.2928f6dea572a6e19d0edfc03e4fde9c41c5d031.de-de.xlf"

# K6 -- "Latest Handback DateTime"
$ws2.Range("K6").Value = "2016-09-01 22:47:32"

# P6 -- "Error Detail"
$ws2.Range("P6").Value = $errorDetail

# Widen the Error Detail column so the message is readable
$ws2.Columns.Item(16).ColumnWidth = 39.17
